# Add a new "test_suite" sheet in front of the existing sheets and
# populate it with the TCID / Runmode table used to drive test skipping.

$wb = $excel.ActiveWorkbook

# Insert the new worksheet before the current first sheet (AddCustomerTest)
# so it becomes the new first tab.
$ws = $wb.Worksheets.Add($wb.Worksheets.Item(1))
$ws.Name = "test_suite"

# Populate the cells in the same order the original authoring tool used,
# so the shared-string table is rebuilt in the same sequence.
$ws.Range("A1").Value = "TCID"
$ws.Range("B1").Value = "Runmode"
$ws.Range("A2").Value = "loginAsBankManagerTest"
$ws.Range("A4").Value = "addCustomerTest"
$ws.Range("B3").Value = "N"
$ws.Range("A3").Value = "openAccountTest"
$ws.Range("B2").Value = "Y"
$ws.Range("B4").Value = "Y"

# Column A sized to fit the longest test-case name.
$ws.Columns.Item(1).ColumnWidth = 19.666666666666668

# Leave the last written cell selected, matching the saved view state.
$ws.Range("B4").Select()
